$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new header cells (row 1) for the new fields ---
# Clone formatting (style "s=3", same as the neighbouring U1/V1 header cells)
# onto the three new header cells before filling in their text/content.
$ws.Range("U1").Copy()
$ws.Range("W1:Y1").PasteSpecial(-4122)

$ws.Range("W1").Value = "الحكم النهائي "
$ws.Range("X1").Value = "تاريخ الحكم النهائي (dd/mm/yyyy) "
$ws.Range("Y1").Value = "رقم الحصر النهائي"

# X1 carries two differently-formatted runs (rich text), matching the other
# "date" headers already present in the sheet (bold headline + smaller
# bold "(dd/mm/yyyy)" hint).
$ws.Range("X1").Characters(1, 20).Font.Bold = $true
$ws.Range("X1").Characters(1, 20).Font.Size = 14
$ws.Range("X1").Characters(21, 13).Font.Bold = $true
$ws.Range("X1").Characters(21, 13).Font.Size = 10

# --- Column sizing for the new columns ---
$ws.Columns.Item(23).ColumnWidth = 13.6667
$ws.Columns.Item(24).ColumnWidth = 22
$ws.Columns.Item(25).ColumnWidth = 13.6667

# --- View: reset scroll position / selection like the saved workbook ---
$ws.Range("B13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
